# Apply updates described in the commit "mejorar apuntes y añadido cnn con mejoras"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: add hours + comment for Jueves 45351 ---
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = "● Leer hasta la página 215 ●Estudio de los ejemplos de predicción stock market con RNN keras y pytorch (apéndice C)"
$ws.Rows.Item(6).RowHeight = 45

# --- Row 8: Sábado 02/03/2024 ---
$ws.Range("C8").Value = "Sábado"
$ws.Range("D8").Value = 45353
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "●Leer hasta la página 280 del libro ●Avanzar apuntes"
$ws.Rows.Item(8).RowHeight = 30

# --- Rows 9-17: fill in day-of-week names and dates ---
$ws.Range("C9").Value = "Domingo"
$ws.Range("D9").Value = 45354

$ws.Range("C10").Value = "Lunes"
$ws.Range("D10").Value = 45355

$ws.Range("C11").Value = "Martes"
$ws.Range("D11").Value = 45356

$ws.Range("C12").Value = "Miércoles"
$ws.Range("D12").Value = 45357

$ws.Range("C13").Value = "Jueves"
$ws.Range("D13").Value = 45358

$ws.Range("C14").Value = "Viernes"
$ws.Range("D14").Value = 45359

$ws.Range("C15").Value = "Sábado"
$ws.Range("D15").Value = 45360

$ws.Range("C16").Value = "Domingo"
$ws.Range("D16").Value = 45361

$ws.Range("C17").Value = "Lunes"
$ws.Range("D17").Value = 45362

# --- Reset the active selection to B1, as in the final workbook state ---
$ws.Range("B1").Select()
